$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 125
$ws_ALC.Range("H125").Value = 975
$ws_ALC.Range("I125").Value = 828.5714
$ws_ALC.Range("J125").Value = 2000
$ws_ALC.Range("K125").Value = 7457.1426
$ws_ALC.Range("L125").Value = 18000
$ws_ALC.Range("M125").Value = -4997.1426
$ws_ALC.Range("N125").Value = -22920

# ALC row 132
$ws_ALC.Range("H132").Value = 1451.3966
$ws_ALC.Range("I132").Value = 1523.0392
$ws_ALC.Range("J132").Value = 929.4286
$ws_ALC.Range("K132").Value = 4569.1176
$ws_ALC.Range("L132").Value = 2788.2858
$ws_ALC.Range("M132").Value = -2039.1176
$ws_ALC.Range("N132").Value = -7848.2858

# ALC row 133
$ws_ALC.Range("H133").Value = 43380
$ws_ALC.Range("J133").Value = 43380
$ws_ALC.Range("L133").Value = 43380
$ws_ALC.Range("N133").Value = -53500

# ALC row 137
$ws_ALC.Range("H137").Value = 871.89655
$ws_ALC.Range("I137").Value = 822.44446
$ws_ALC.Range("J137").Value = 894.15
$ws_ALC.Range("K137").Value = 2467.33338
$ws_ALC.Range("L137").Value = 2682.45
$ws_ALC.Range("M137").Value = 82.66661999999997
$ws_ALC.Range("N137").Value = -7782.45

# ALC row 138
$ws_ALC.Range("H138").Value = 2022.39
$ws_ALC.Range("I138").Value = 892.7805
$ws_ALC.Range("J138").Value = 2807.3728
$ws_ALC.Range("K138").Value = 2678.3415
$ws_ALC.Range("L138").Value = 8422.118399999999
$ws_ALC.Range("M138").Value = 2461.6585
$ws_ALC.Range("N138").Value = -18702.1184

# ALC row 141
$ws_ALC.Range("H141").Value = 2858.389
$ws_ALC.Range("I141").Value = 1189
$ws_ALC.Range("J141").Value = 8701.25
$ws_ALC.Range("K141").Value = 3567
$ws_ALC.Range("L141").Value = 26103.75
$ws_ALC.Range("M141").Value = 1613
$ws_ALC.Range("N141").Value = -36463.75

# ARM row 32
$ws_ARM.Range("H32").Value = 19969.186
$ws_ARM.Range("I32").Value = 18977.393
$ws_ARM.Range("J32").Value = 24322.055
$ws_ARM.Range("K32").Value = 18977.393
$ws_ARM.Range("L32").Value = 24322.055
$ws_ARM.Range("M32").Value = -18690.393
$ws_ARM.Range("N32").Value = -24896.055

# ARM row 61
$ws_ARM.Range("H61").Value = 866.72974
$ws_ARM.Range("I61").Value = 730.59375
$ws_ARM.Range("J61").Value = 1738
$ws_ARM.Range("K61").Value = 730.59375
$ws_ARM.Range("L61").Value = 1738
$ws_ARM.Range("M61").Value = -518.59375
$ws_ARM.Range("N61").Value = -2162

# ARM row 132
$ws_ARM.Range("H132").Value = 1152.5128
$ws_ARM.Range("I132").Value = 1022.4
$ws_ARM.Range("J132").Value = 1586.2222
$ws_ARM.Range("K132").Value = 3067.2
$ws_ARM.Range("L132").Value = 4758.6666
$ws_ARM.Range("M132").Value = -537.1999999999998
$ws_ARM.Range("N132").Value = -9818.6666

# ARM row 136
$ws_ARM.Range("H136").Value = 866.72974
$ws_ARM.Range("I136").Value = 730.59375
$ws_ARM.Range("J136").Value = 1738
$ws_ARM.Range("K136").Value = 2191.78125
$ws_ARM.Range("L136").Value = 5214
$ws_ARM.Range("M136").Value = 358.21875
$ws_ARM.Range("N136").Value = -10314

# BSM row 44
$ws_BSM.Range("H44").Value = 0
$ws_BSM.Range("J44").Value = 0
$ws_BSM.Range("L44").Value = 0
$ws_BSM.Range("N44").ClearContents()

# BSM row 134
$ws_BSM.Range("H134").Value = 20453.389
$ws_BSM.Range("I134").Value = 1636.5333
$ws_BSM.Range("J134").Value = 114537.664
$ws_BSM.Range("K134").Value = 4909.5999
$ws_BSM.Range("L134").Value = 343612.992
$ws_BSM.Range("M134").Value = -2374.5999
$ws_BSM.Range("N134").Value = -348682.992

# CRP row 31
$ws_CRP.Range("H31").Value = 2688.8718
$ws_CRP.Range("I31").Value = 2598.9033
$ws_CRP.Range("J31").Value = 3037.5
$ws_CRP.Range("K31").Value = 2598.9033
$ws_CRP.Range("L31").Value = 3037.5
$ws_CRP.Range("M31").Value = -2303.9033
$ws_CRP.Range("N31").Value = -3627.5

# CRP row 34
$ws_CRP.Range("H34").Value = 2688.8718
$ws_CRP.Range("I34").Value = 2598.9033
$ws_CRP.Range("J34").Value = 3037.5
$ws_CRP.Range("K34").Value = 2598.9033
$ws_CRP.Range("L34").Value = 3037.5
$ws_CRP.Range("M34").Value = -2396.9033
$ws_CRP.Range("N34").Value = -3441.5

# CRP row 53
$ws_CRP.Range("H53").Value = 23000
$ws_CRP.Range("J53").Value = 23000
$ws_CRP.Range("L53").Value = 23000
$ws_CRP.Range("N53").Value = -24214

# CRP row 58
$ws_CRP.Range("H58").Value = 3131.0833
$ws_CRP.Range("I58").Value = 1037.0303
$ws_CRP.Range("J58").Value = 7738
$ws_CRP.Range("K58").Value = 1037.0303
$ws_CRP.Range("L58").Value = 7738
$ws_CRP.Range("M58").Value = -834.0302999999999
$ws_CRP.Range("N58").Value = -8144

# CRP row 62
$ws_CRP.Range("H62").Value = 5474.6
$ws_CRP.Range("I62").Value = 5559.5293
$ws_CRP.Range("J62").Value = 4993.3335
$ws_CRP.Range("K62").Value = 5559.5293
$ws_CRP.Range("L62").Value = 4993.3335
$ws_CRP.Range("M62").Value = -4935.5293
$ws_CRP.Range("N62").Value = -6241.3335

# CRP row 65
$ws_CRP.Range("H65").Value = 5474.6
$ws_CRP.Range("I65").Value = 5559.5293
$ws_CRP.Range("J65").Value = 4993.3335
$ws_CRP.Range("K65").Value = 27797.6465
$ws_CRP.Range("L65").Value = 24966.6675
$ws_CRP.Range("M65").Value = -24677.6465
$ws_CRP.Range("N65").Value = -31206.6675

# CRP row 132
$ws_CRP.Range("H132").Value = 1842.4354
$ws_CRP.Range("I132").Value = 1175.8529
$ws_CRP.Range("J132").Value = 2651.8572
$ws_CRP.Range("K132").Value = 3527.5587
$ws_CRP.Range("L132").Value = 7955.571599999999
$ws_CRP.Range("M132").Value = -997.5587000000005
$ws_CRP.Range("N132").Value = -13015.5716

# CRP row 134
$ws_CRP.Range("H134").Value = 1273.2069
$ws_CRP.Range("I134").Value = 1220.8
$ws_CRP.Range("K134").Value = 3662.4
$ws_CRP.Range("M134").Value = -1127.4

# CRP row 136
$ws_CRP.Range("H136").Value = 3131.0833
$ws_CRP.Range("I136").Value = 1037.0303
$ws_CRP.Range("J136").Value = 7738
$ws_CRP.Range("K136").Value = 3111.0909
$ws_CRP.Range("L136").Value = 23214
$ws_CRP.Range("M136").Value = -561.0908999999997
$ws_CRP.Range("N136").Value = -28314

# CUL row 5
$ws_CUL.Range("H5").Value = 849.3214
$ws_CUL.Range("I5").Value = 836.9167
$ws_CUL.Range("J5").Value = 923.75
$ws_CUL.Range("K5").Value = 2510.7501
$ws_CUL.Range("L5").Value = 2771.25
$ws_CUL.Range("M5").Value = -2398.7501
$ws_CUL.Range("N5").Value = -2995.25

# CUL row 74
$ws_CUL.Range("H74").Value = 4998.5713
$ws_CUL.Range("J74").Value = 4998.5713
$ws_CUL.Range("L74").Value = 14995.7139
$ws_CUL.Range("N74").Value = -17117.7139

# CUL row 77
$ws_CUL.Range("H77").Value = 4998.5713
$ws_CUL.Range("J77").Value = 4998.5713
$ws_CUL.Range("L77").Value = 44987.14169999999
$ws_CUL.Range("N77").Value = -55595.14169999999

# CUL row 131
$ws_CUL.Range("H131").Value = 30395.941
$ws_CUL.Range("J131").Value = 18180.568
$ws_CUL.Range("L131").Value = 54541.704
$ws_CUL.Range("N131").Value = -64621.704

# CUL row 135
$ws_CUL.Range("H135").Value = 849.3214
$ws_CUL.Range("I135").Value = 836.9167
$ws_CUL.Range("J135").Value = 923.75
$ws_CUL.Range("K135").Value = 7532.2503
$ws_CUL.Range("L135").Value = 8313.75
$ws_CUL.Range("M135").Value = -4997.2503
$ws_CUL.Range("N135").Value = -13383.75

# GSM row 140
$ws_GSM.Range("H140").Value = 71000
$ws_GSM.Range("J140").Value = 71000
$ws_GSM.Range("L140").Value = 71000
$ws_GSM.Range("N140").Value = -81360

# LTW row 22
$ws_LTW.Range("H22").Value = 1049.4584
$ws_LTW.Range("I22").Value = 374.06668
$ws_LTW.Range("J22").Value = 1356.4546
$ws_LTW.Range("K22").Value = 374.06668
$ws_LTW.Range("L22").Value = 1356.4546
$ws_LTW.Range("M22").Value = -79.06668000000002
$ws_LTW.Range("N22").Value = -1946.4546

# LTW row 27
$ws_LTW.Range("H27").Value = 1049.4584
$ws_LTW.Range("I27").Value = 374.06668
$ws_LTW.Range("J27").Value = 1356.4546
$ws_LTW.Range("K27").Value = 374.06668
$ws_LTW.Range("L27").Value = 1356.4546
$ws_LTW.Range("M27").Value = -267.06668
$ws_LTW.Range("N27").Value = -1570.4546

# LTW row 61
$ws_LTW.Range("H61").Value = 4582.645
$ws_LTW.Range("I61").Value = 5144.6924
$ws_LTW.Range("J61").Value = 1660
$ws_LTW.Range("K61").Value = 5144.6924
$ws_LTW.Range("L61").Value = 1660
$ws_LTW.Range("M61").Value = -4942.6924
$ws_LTW.Range("N61").Value = -2064

# LTW row 82
$ws_LTW.Range("H82").Value = 1493.862
$ws_LTW.Range("I82").Value = 1710.75
$ws_LTW.Range("J82").Value = 1226.9231
$ws_LTW.Range("K82").Value = 1710.75
$ws_LTW.Range("L82").Value = 1226.9231
$ws_LTW.Range("M82").Value = -1349.75
$ws_LTW.Range("N82").Value = -1948.9231

# LTW row 85
$ws_LTW.Range("H85").Value = 1493.862
$ws_LTW.Range("I85").Value = 1710.75
$ws_LTW.Range("J85").Value = 1226.9231
$ws_LTW.Range("K85").Value = 1710.75
$ws_LTW.Range("L85").Value = 1226.9231
$ws_LTW.Range("M85").Value = -462.75
$ws_LTW.Range("N85").Value = -3722.9231

# LTW row 113
$ws_LTW.Range("H113").Value = 4582.645
$ws_LTW.Range("I113").Value = 5144.6924
$ws_LTW.Range("J113").Value = 1660
$ws_LTW.Range("K113").Value = 5144.6924
$ws_LTW.Range("L113").Value = 1660
$ws_LTW.Range("M113").Value = -2974.6924
$ws_LTW.Range("N113").Value = -6000

# LTW row 132
$ws_LTW.Range("H132").Value = 2464.1707
$ws_LTW.Range("I132").Value = 2379.5356
$ws_LTW.Range("J132").Value = 2646.4614
$ws_LTW.Range("K132").Value = 7138.6068
$ws_LTW.Range("L132").Value = 7939.3842
$ws_LTW.Range("M132").Value = -4608.6068
$ws_LTW.Range("N132").Value = -12999.3842

# LTW row 136
$ws_LTW.Range("H136").Value = 1648.1587
$ws_LTW.Range("I136").Value = 929.8868
$ws_LTW.Range("J136").Value = 5455
$ws_LTW.Range("K136").Value = 2789.6604
$ws_LTW.Range("L136").Value = 16365
$ws_LTW.Range("M136").Value = -239.6603999999998
$ws_LTW.Range("N136").Value = -21465

# WVR row 132
$ws_WVR.Range("H132").Value = 792.0599999999999
$ws_WVR.Range("I132").Value = 632.8108
$ws_WVR.Range("J132").Value = 1245.3077
$ws_WVR.Range("K132").Value = 1898.4324
$ws_WVR.Range("L132").Value = 3735.9231
$ws_WVR.Range("M132").Value = 631.5676000000001
$ws_WVR.Range("N132").Value = -8795.9231
